$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Descuento" column (L) header
$ws.Range("L1").Value = "Descuento"
$ws.Range("L1").Font.Size = 12

# Discount value for existing product (row 2)
$ws.Range("L2").Value = 5
$ws.Range("L2").Font.Size = 12

# New product row (row 3): "TG pasta gypsum galon (masilla)"
$ws.Range("B3").Value = 13002
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "TG pasta gypsum galon (masilla)"
$ws.Range("D3").Font.Size = 12
$ws.Range("E3").Value = 3876.11
$ws.Range("F3").Value = "galon"
$ws.Range("F3").Font.Size = 12
$ws.Range("G3").Value = "t"
$ws.Range("G3").Font.Size = 12
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 13
$ws.Range("L3").Value = 10
$ws.Range("L3").Font.Size = 12

# Restore selection similar to the authored workbook
$ws.Range("E17").Select() | Out-Null
